$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pharmacy name title from the merged header cell (D2:N2)
$ws.Range("D2").Value = ""

# Update the printed timestamp in A10 to reflect the new save time
$ws.Range("A10").Value = "Wednesday, 17 September, 2025 10:59 PM"

# Refresh row heights to match Excel's autofit after the header text change
$ws.Rows.Item(2).RowHeight = 25.5
$ws.Rows.Item(3).RowHeight = 35.25
$ws.Rows.Item(5).RowHeight = 0.75
$ws.Rows.Item(6).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 26.25
